$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$arr = New-Object 'object[,]' 100,2
$arr[0,0] = -0.6436936677321438
$arr[0,1] = 0.6048615524379573
$arr[1,0] = 0.8481987922759403
$arr[1,1] = -0.06183975955357971
$arr[2,0] = 0.1362243115762983
$arr[2,1] = 0.6668242860948045
$arr[3,0] = -2.866014234017723
$arr[3,1] = 1.427542058963176
$arr[4,0] = -0.4341914718821441
$arr[4,1] = -1.875270068897947
$arr[5,0] = -1.448678003288501
$arr[5,1] = -0.7410904468199772
$arr[6,0] = -2.670841336409842
$arr[6,1] = 0.3023040868638705
$arr[7,0] = 0.253343792392148
$arr[7,1] = 0.3395936605247907
$arr[8,0] = 1.124378825728025
$arr[8,1] = 0.5815546117344512
$arr[9,0] = -0.3328171023833046
$arr[9,1] = -0.734141001411168
$arr[10,0] = 0.9715251082647074
$arr[10,1] = -0.5451462484222711
$arr[11,0] = -0.3151471034450695
$arr[11,1] = 0.9682226724522393
$arr[12,0] = -0.4595775492837924
$arr[12,1] = 0.8289935256924429
$arr[13,0] = 0.6593929085579789
$arr[13,1] = -2.243730152265649
$arr[14,0] = -0.07791176058193165
$arr[14,1] = -0.6495241628031473
$arr[15,0] = 0.1686142624598296
$arr[15,1] = 0.191724201537998
$arr[16,0] = 2.382645102519338
$arr[16,1] = -0.4582057661330416
$arr[17,0] = 0.550240542470282
$arr[17,1] = 1.436781058176872
$arr[18,0] = -1.32313431129281
$arr[18,1] = -2.230893663151085
$arr[19,0] = 0.1201625389229685
$arr[19,1] = 0.6169262894177493
$arr[20,0] = -1.238371232078775
$arr[20,1] = -0.3955137710785661
$arr[21,0] = -1.794652629800652
$arr[21,1] = 2.468249722939703
$arr[22,0] = -0.0836666869437588
$arr[22,1] = -0.9638709116197721
$arr[23,0] = -1.479878989264926
$arr[23,1] = 0.3508202005104995
$arr[24,0] = 1.482014008149077
$arr[24,1] = 0.7554382049861887
$arr[25,0] = -0.4152403562579223
$arr[25,1] = 1.148642674641205
$arr[26,0] = 1.123590593476252
$arr[26,1] = -1.060910880374953
$arr[27,0] = 0.1594692143239325
$arr[27,1] = 1.221488544430398
$arr[28,0] = -0.2933732025751061
$arr[28,1] = -0.214152447243255
$arr[29,0] = -0.8855557045211763
$arr[29,1] = 1.020915229079611
$arr[30,0] = 0.0620496432929693
$arr[30,1] = 0.3544165865739892
$arr[31,0] = -1.3702988547915
$arr[31,1] = -0.6461119001904715
$arr[32,0] = -1.754582786087739
$arr[32,1] = -0.3991244150749291
$arr[33,0] = 0.5664072676744579
$arr[33,1] = -1.480385013511401
$arr[34,0] = 1.445166997699438
$arr[34,1] = -1.639949060347753
$arr[35,0] = 0.4459971448962949
$arr[35,1] = -0.5379657233611839
$arr[36,0] = -1.565073502983542
$arr[36,1] = -0.6950834099672214
$arr[37,0] = 0.2463880731422886
$arr[37,1] = -0.2477069577402639
$arr[38,0] = -0.3507242913887475
$arr[38,1] = -1.44051868664003
$arr[39,0] = -0.6514460131826214
$arr[39,1] = -0.6292782376388897
$arr[40,0] = -0.5464757699054879
$arr[40,1] = -0.1714765588734759
$arr[41,0] = 0.3231629392285399
$arr[41,1] = 0.6534159802800087
$arr[42,0] = -0.2915282051470537
$arr[42,1] = 0.6967750831528265
$arr[43,0] = 0.6151820722118405
$arr[43,1] = 1.091175472278164
$arr[44,0] = 2.588000478193001
$arr[44,1] = -0.04589524596991784
$arr[45,0] = 1.370749451107059
$arr[45,1] = -1.178660215603311
$arr[46,0] = 0.9393252668458343
$arr[46,1] = -1.355580438879415
$arr[47,0] = 0.3006316064031579
$arr[47,1] = 0.771349844936364
$arr[48,0] = 0.05815274417296564
$arr[48,1] = 2.246041315098033
$arr[49,0] = 0.5267149621358643
$arr[49,1] = 0.8233417008525109
$arr[50,0] = -0.7123731816434228
$arr[50,1] = 0.276469818975218
$arr[51,0] = 0.8581973449076185
$arr[51,1] = 0.2413355351084589
$arr[52,0] = -1.145956704531524
$arr[52,1] = -1.226158624623136
$arr[53,0] = 0.5361686707597
$arr[53,1] = -0.2499597238895034
$arr[54,0] = 0.7660601344380447
$arr[54,1] = 0.6634042422292832
$arr[55,0] = -0.9432762958309731
$arr[55,1] = 0.2823738417461684
$arr[56,0] = -0.837221143630809
$arr[56,1] = -0.8561477056658795
$arr[57,0] = -0.8351698402828253
$arr[57,1] = 0.8263068678054736
$arr[58,0] = 0.1419408932370512
$arr[58,1] = -1.495064388596024
$arr[59,0] = 0.3022392408561329
$arr[59,1] = -0.3825669144857662
$arr[60,0] = -0.4127490461341597
$arr[60,1] = -0.6458329774761561
$arr[61,0] = -2.085492571584192
$arr[61,1] = 0.3716365098571095
$arr[62,0] = 0.1315193023423922
$arr[62,1] = 1.222406872154715
$arr[63,0] = -0.2386606521516989
$arr[63,1] = 0.5866952124266063
$arr[64,0] = 0.2617963838242409
$arr[64,1] = 0.06977936609247397
$arr[65,0] = -0.6241851056429294
$arr[65,1] = -1.037653362277436
$arr[66,0] = -0.4539420747994777
$arr[66,1] = 1.30576054228318
$arr[67,0] = -0.404701238277281
$arr[67,1] = 1.075223412767373
$arr[68,0] = 1.235252177085605
$arr[68,1] = -0.8031512192257261
$arr[69,0] = 0.09702201463857606
$arr[69,1] = 1.453845934993659
$arr[70,0] = -0.09487154581914876
$arr[70,1] = 1.814408098157268
$arr[71,0] = 0.2271801367917244
$arr[71,1] = -1.115393088223188
$arr[72,0] = 0.9328267110342523
$arr[72,1] = 2.168521318056309
$arr[73,0] = -0.09315705105604782
$arr[73,1] = -0.345436021703852
$arr[74,0] = 1.878318642470408
$arr[74,1] = 1.116318477641528
$arr[75,0] = 0.1932647763661189
$arr[75,1] = -0.5724231192636764
$arr[76,0] = 0.1912617951639147
$arr[76,1] = -0.6834481631749896
$arr[77,0] = -0.2693225846762528
$arr[77,1] = 0.4378468894079781
$arr[78,0] = 0.5623434878191472
$arr[78,1] = 0.3404363215159713
$arr[79,0] = -0.8391119319658223
$arr[79,1] = -1.456319730170735
$arr[80,0] = -1.948869760072637
$arr[80,1] = -0.1042592799032529
$arr[81,0] = -1.44804662971959
$arr[81,1] = -0.02036165746055315
$arr[82,0] = 0.06947472423285476
$arr[82,1] = 0.3447790295087322
$arr[83,0] = 0.03229711418411784
$arr[83,1] = -1.180213158933638
$arr[84,0] = -0.03350886978972536
$arr[84,1] = 0.03624222666411998
$arr[85,0] = 0.7344473818902533
$arr[85,1] = -0.510248074893689
$arr[86,0] = -1.092627473367835
$arr[86,1] = -0.7863097810137865
$arr[87,0] = 0.5658131615087325
$arr[87,1] = 1.195582748791698
$arr[88,0] = 0.755411988513544
$arr[88,1] = 0.1232495374633156
$arr[89,0] = -0.678145361754754
$arr[89,1] = -1.102240414597044
$arr[90,0] = -1.916282923517614
$arr[90,1] = 0.7102227163283419
$arr[91,0] = -0.8248177828102289
$arr[91,1] = 0.3554956757573015
$arr[92,0] = 1.228226554465794
$arr[92,1] = 0.03827843410827805
$arr[93,0] = -0.5431044208630181
$arr[93,1] = -1.12796578461849
$arr[94,0] = 0.3684630369346275
$arr[94,1] = 1.277668742749809
$arr[95,0] = 1.287616710102915
$arr[95,1] = -0.5011707599454829
$arr[96,0] = -0.07176661507610869
$arr[96,1] = 1.301262065439911
$arr[97,0] = 0.7745833125790254
$arr[97,1] = -0.7349295588644752
$arr[98,0] = 0.6957844459985409
$arr[98,1] = -1.334993941807255
$arr[99,0] = -0.2988073625544497
$arr[99,1] = -0.4763227695763352
$ws.Range("B2:C101").Value = $arr
Write-Host "done"